$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.548.38'
$ws.Range('E2').Value = '  +5.19%  '
$ws.Range('D3').Value = '2.254.09'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'232.01"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').Value = "'0.638"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.41%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').Value = "'59.44"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('D11').Value = "'0.0904"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.62%  '
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '2.587.73'
$ws.Range('E13').Value = '  +4.37%  '
$ws.Range('D14').Value = "'16.25"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.96%  '
$ws.Range('D15').Value = "'22.65"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').Value = "'0.834"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('D17').Value = "'5.68"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').Value = '2.262.50'
$ws.Range('E18').Value = '  +5.02%  '
$ws.Range('D19').Value = '41.429.37'
$ws.Range('E19').Value = '  +5.00%  '
$ws.Range('D20').Value = "'73.87"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.77%  '
$ws.Range('D21').Value = '0.0₃0919'
$ws.Range('E21').Value = '  +7.89%  '
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('D23').Value = "'251.58"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +9.43%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.40"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'2.35"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('D27').Value = "'9.78"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('D28').Value = "'173.42"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('E29').Value = '  +3.23%  '
$ws.Range('D30').Value = "'20.52"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.22%  '
$ws.Range('E31').Value = '  +2.48%  '
$ws.Range('D32').Value = "'2.81"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +8.21%  '
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('D34').Value = "'5.06"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.89%  '
$ws.Range('D35').Value = "'4.76"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.48%  '
$ws.Range('D36').Value = "'0.0636"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.18%  '
$ws.Range('D37').Value = "'7.04"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('D38').Value = "'3.85"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.89%  '
$ws.Range('D39').Value = "'2.47"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').Value = "'0.000269"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +72.34%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = "'4.95"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +14.36%  '
$ws.Range('D43').Value = "'0.0241"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.13%  '
$ws.Range('D44').Value = "'8.83"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +15.00%  '
$ws.Range('D45').Value = "'102.88"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').Value = "'17.88"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').Value = "'1.23"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.93%  '
$ws.Range('D48').Value = '1.511.57'
$ws.Range('D49').Value = "'0.0946"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('D51').Value = "'2.79"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.25%  '
